$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: TOTAL (5-17 y.o.) / All population groups ---
$ws.Range("C2").Value = 68663199
$ws.Range("D2").Value = 65.59999999999999
$ws.Range("E2").Value = 45060900
$ws.Range("F2").Value = 33.4
$ws.Range("G2").Value = 22928793
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 665846
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7661
$ws.Range("L2").Value = 34.4
$ws.Range("M2").Value = 23602299

# --- Row 3: ocap -> non_pdi ---
$ws.Range("A3").Value = "non_pdi (5-17 y.o.)"
$ws.Range("B3").Value = "non_pdi"
$ws.Range("C3").Value = 48523504
$ws.Range("D3").Value = 74.2
$ws.Range("E3").Value = 35998252
$ws.Range("F3").Value = 24.8
$ws.Range("G3").Value = 12049923
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 468923
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6406
$ws.Range("L3").Value = 25.8
$ws.Range("M3").Value = 12525253

# --- Row 4: idp -> pdi ---
$ws.Range("A4").Value = "pdi (5-17 y.o.)"
$ws.Range("B4").Value = "pdi"
$ws.Range("C4").Value = 20139695
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 9062649
$ws.Range("F4").Value = 54
$ws.Range("G4").Value = 10878870
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 196922
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1254
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 11077047

# --- Remove rows 5 and 6 (ret, ndsp) entirely ---
$ws.Range("A5:M6").EntireRow.Delete()
